$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.025.68'
$ws.Range("E2").Value = '  +2.66%  '
$ws.Range("D3").Value = '2.234.12'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''294.42'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = '''86.48'
$ws.Range("E6").Value = '  +8.35%  '
$ws.Range("D7").Value = '''0.515'
$ws.Range("E7").Value = '  +2.44%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.472'
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("D10").Value = '''31.12'
$ws.Range("E10").Value = '  +11.29%  '
$ws.Range("D11").Value = '''0.0792'
$ws.Range("E11").Value = '  +2.53%  '
$ws.Range("D12").Value = '''47.08'
$ws.Range("E12").Value = '  +3.12%  '
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").Value = '''6.46'
$ws.Range("E14").Value = '  +6.19%  '
$ws.Range("D15").Value = '2.582.51'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = '''14.15'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").Value = '2.234.82'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '''0.729'
$ws.Range("E18").Value = '  +2.99%  '
$ws.Range("D19").Value = '39.963.56'
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  +3.65%  '
$ws.Range("D21").Value = '''5.80'
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").Value = '''10.82'
$ws.Range("E22").Value = '  +11.02%  '
$ws.Range("D23").Value = '''65.46'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").Value = '''235.18'
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").Value = '''2.46'
$ws.Range("E26").Value = '  +3.85%  '
$ws.Range("D27").Value = '''1.84'
$ws.Range("E27").Value = '  +5.89%  '
$ws.Range("D28").Value = '''22.79'
$ws.Range("E28").Value = '  +2.62%  '
$ws.Range("E29").Value = '  +3.19%  '
$ws.Range("D30").Value = '''9.22'
$ws.Range("E30").Value = '  +3.71%  '
$ws.Range("D31").Value = '''33.30'
$ws.Range("E31").Value = '  +7.24%  '
$ws.Range("D32").Value = '''152.36'
$ws.Range("E32").Value = '  +2.56%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '''4.88'
$ws.Range("E34").Value = '  +3.41%  '
$ws.Range("D35").Value = '''0.0716'
$ws.Range("E35").Value = '  +4.72%  '
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("D37").Value = '''16.37'
$ws.Range("E37").Value = '  +14.29%  '
$ws.Range("D38").Value = '''0.111'
$ws.Range("E38").Value = '  +3.19%  '
$ws.Range("E39").Value = '  +5.55%  '
$ws.Range("D40").Value = '''2.71'
$ws.Range("E40").Value = '  +3.07%  '
$ws.Range("D41").Value = '''1.69'
$ws.Range("E41").Value = '  +7.00%  '
$ws.Range("D42").Value = '''3.84'
$ws.Range("E42").Value = '  +7.09%  '
$ws.Range("D43").Value = '2.042.15'
$ws.Range("E43").Value = '  +7.46%  '
$ws.Range("D44").Value = '''2.24'
$ws.Range("E44").Value = '  +7.95%  '
$ws.Range("D45").Value = '''0.0270'
$ws.Range("E45").Value = '  +6.71%  '
$ws.Range("D46").Value = '''9.97'
$ws.Range("E46").Value = '  +13.06%  '
$ws.Range("D47").Value = '''16.33'
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("D48").Value = '''2.56'
$ws.Range("E48").Value = '  +2.55%  '
$ws.Range("D49").Value = '2.452.88'
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("D50").Value = '''70.87'
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("D51").Value = '''1.45'
$ws.Range("E51").Value = '  +14.31%  '
